# Matej.docx edit:
#   - paragraph "dorobiť funkciu get_user_id... - DONE" loses its trailing
#     " - " (hyphen) in favour of " – " (en dash) and the _GoBack bookmark
#     that used to wrap "DONE" is dropped from this paragraph.
#   - a brand-new bullet is inserted right after it:
#       "Prerobit priority ciselne (strojove) na pisane (ludske) - DONE"
#     (same ListParagraph / numbering / yellow highlight as its neighbours,
#     "DONE" bold, and the _GoBack bookmark now sitting empty in the
#     middle of "pisane", right after "pis").

$d = $word.ActiveDocument

# 1) Locate the "dorobiť funkciu get_user_id..." bullet by its distinctive
#    text instead of a hard-coded paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*get_user_id*") {
        $target = $cand
        break
    }
}

# 2) Swap the " - " right before "DONE" for " – " (en dash). Search is
#    scoped to this single paragraph so it cannot touch the bold "DONE"
#    run (kept as a separate run) or bleed into neighbouring bullets.
$searchRng = $target.Range
$searchRng.Find.Execute(" - ", $false, $false, $false, $false, $false, `
    $true, 1, $false, " " + [char]0x2013 + " ", 2)

# Re-resolve the (now en-dash) paragraph and append a new paragraph right
# after it, inheriting its ListParagraph style / numbering / highlight.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*get_user_id*") {
        $target = $cand
        break
    }
}
$newParaIndex = $target.Index + 1
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($newParaIndex)
# The freshly split paragraph mark inherited "Bold" from the "DONE" run it
# was split off from - clear it before typing so the new sentence starts
# out un-bold (only the trailing "DONE" should be bold).
$newPara.Range.Font.Bold = 0
$newPara.Range.InsertAfter( `
    "Prerobit priority ciselne (strojove) na pisane (ludske) - DONE")

# Re-resolve the paragraph and bold just the trailing "DONE".
$newPara = $d.Paragraphs.Item($newParaIndex)
$doneEnd = $newPara.Range.End - 1
$doneStart = $doneEnd - 4
$doneRng = $d.Range($doneStart, $doneEnd)
$doneRng.Font.Bold = 1

# Move the "_GoBack" bookmark: it used to wrap "DONE" in the previous
# paragraph; now it is an empty bookmark sitting right after "pis" in
# "pisane" on the new line ("Prerobit priority ciselne (strojove) na
# pis|ane (ludske) - DONE").
$newPara = $d.Paragraphs.Item($newParaIndex)
$bmPos = $newPara.Range.Start + [string]"Prerobit priority ciselne (strojove) na pis".Length
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
